$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D28").Value = "QMIX: Monotonic Value Function Factorisation for Deep Multi-Agent Reinforcement Learning 리뷰"
$ws.Range("E28").Value = "https://ropiens.tistory.com/112"

$ws.Range("D36").Value = "Introduction to Steel Surface Defect Detection"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/321"

$ws.Range("D44").Value = "AI 알고리즘 경량화"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/79"

$ws.Range("D51").Value = "[MariaDB] 스토어드 프로시저 사용법"
$ws.Range("E51").Value = "https://bskyvision.com/1187"
